# OpenData_Slovakia_Covid_VaccinesDelivery.xlsx - "Add files via upload"
# Corrected delivery dates/quantities for the last few J&J (col F) delivery rows
# and updated status/status-update info for rows 73-74.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 73: date pulled back 3 days, dose count revised up, status moved from
# "Assumption" to "Confirmed", and a concrete status-update date recorded.
$ws.Range("A73").Value = 44330
$ws.Range("C73").Value = 12000
$ws.Range("D73").Value = "Confirmed"
$ws.Range("E73").Value = 44320

# Row 74: same kind of update as row 73.
$ws.Range("A74").Value = 44337
$ws.Range("C74").Value = 12000
$ws.Range("D74").Value = "Confirmed"
$ws.Range("E74").Value = 44328

# Rows 75-79: only the delivery date shifts back by 3 days.
$ws.Range("A75").Value = 44344
$ws.Range("A76").Value = 44351
$ws.Range("A77").Value = 44358
$ws.Range("A78").Value = 44365
$ws.Range("A79").Value = 44372

# Refresh the sheet view: scroll down so row 61 is at the top and select the
# full used range (A1:F79), matching the view state saved with this edit.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 61
$ws.Range("A1:F79").Select()
